$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains text formatting so numeric-looking
# strings (e.g. "1.0000", "29.841.27") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.841.27'
$ws.Range("D3").Value = '1.898.81'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '0.7782'
$ws.Range("E5").Value = '  +5.38%  '
$ws.Range("D6").Value = '240.44'
$ws.Range("E6").Value = '  -1.48%  '
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.3067'
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("E9").Value = '  -5.02%  '
$ws.Range("D10").Value = '0.06861'
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").Value = '0.07990'
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("D12").Value = '1.914.14'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '0.7386'
$ws.Range("E13").Value = '  -5.17%  '
$ws.Range("D14").Value = '5.182'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '29.860.70'
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = '13.79'
$ws.Range("E17").Value = '  -4.43%  '
$ws.Range("D18").Value = '5.896'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("D19").Value = '244.84'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = '0.000007718'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D22").Value = '2.150.66'
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '6.935'
$ws.Range("E24").Value = '  -3.74%  '
$ws.Range("D25").Value = '9.283'
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("D26").Value = '166.85'
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").Value = '18.76'
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").Value = '0.1307'
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("D29").Value = '2.029'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").Value = '1.389'
$ws.Range("E30").Value = '  +2.51%  '
$ws.Range("D31").Value = '1.511'
$ws.Range("E31").Value = '  -2.39%  '
$ws.Range("D32").Value = '4.280'
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").Value = '4.070'
$ws.Range("E33").Value = '  -0.87%  '
$ws.Range("D34").Value = '0.05254'
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("D36").Value = '0.7288'
$ws.Range("E36").Value = '  -3.24%  '
$ws.Range("D37").Value = '2.727'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  -1.91%  '
$ws.Range("D39").Value = '2.781'
$ws.Range("E39").Value = '  -0.89%  '
$ws.Range("D40").Value = '6.193'
$ws.Range("E40").Value = '  -3.05%  '
$ws.Range("D41").Value = '0.4423'
$ws.Range("E41").Value = '  -2.05%  '
$ws.Range("D42").Value = '72.24'
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '0.8384'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '1.882'
$ws.Range("E45").Value = '  -4.38%  '
$ws.Range("D46").Value = '7.594'
$ws.Range("E46").Value = '  -3.39%  '
$ws.Range("D47").Value = '100.37'
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("D48").Value = '9.787'
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("D49").Value = '2.065.13'
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").Value = '36.18'
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("D51").Value = '930.22'
$ws.Range("E51").Value = '  -0.78%  '
